$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume/1h) hold numeric-looking / percentage-looking
# text values that must stay as plain text (matching the original inlineStr cells),
# so we force a temporary "@" (Text) number format before assigning the value, then
# restore the default "Normal" style so the cell format matches the original (no style).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '296.71'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.02%'
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '41.74'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.34%'
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.034'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.09%'
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07555'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.60%'
$ws.Range('E5').Style = 'Normal'

$ws.Range('B6').Value = 'GateToken'

$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.389'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2.46%'
$ws.Range('E6').Style = 'Normal'

$ws.Range('B7').Value = 'FTXToken'

$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.599'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.25%'
$ws.Range('E7').Style = 'Normal'

$ws.Range('B8').Value = 'MXToken'

$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9289'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.91%'
$ws.Range('E8').Style = 'Normal'

$ws.Range('B9').Value = 'BTSEToken'

$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.408'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2.13%'
$ws.Range('E9').Style = 'Normal'

$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1197'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '3.22%'
$ws.Range('E10').Style = 'Normal'

$ws.Range('B11').Value = 'WazirX'

$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1840'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '6.42%'
$ws.Range('E11').Style = 'Normal'

$ws.Range('B12').Value = 'MandalaExchangeToken'

$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09038'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '4.33%'
$ws.Range('E12').Style = 'Normal'

$ws.Range('B13').Value = 'BitrueCoin'

$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03974'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-5.09%'
$ws.Range('E13').Style = 'Normal'

$ws.Range('B14').Value = 'BitMartToken'

$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.1053'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.10%'
$ws.Range('E14').Style = 'Normal'

$ws.Range('B15').Value = 'BitForexToken'

$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001292'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '2.18%'
$ws.Range('E15').Style = 'Normal'

$ws.Range('B16').Value = 'TigerCash'

$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005934'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '2.14%'
$ws.Range('E16').Style = 'Normal'

$ws.Range('B17').Value = 'LEO'

$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.353'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-1.48%'
$ws.Range('E17').Style = 'Normal'

$ws.Range('B18').Value = 'BitpandaEcosystemToken'

$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.3320'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '1.31%'
$ws.Range('E18').Style = 'Normal'

$ws.Range('B19').Value = 'MCDex'

$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.906'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.72%'
$ws.Range('E19').Style = 'Normal'

$ws.Range('B20').Value = 'ProBitToken'

$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1420'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '2.78%'
$ws.Range('E20').Style = 'Normal'

$ws.Range('B21').Value = 'ZBToken'

$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.3000'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '4.01%'
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04056'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '5.14%'
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.37%'
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.003933'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '3.64%'
$ws.Range('E24').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.07%'
$ws.Range('E26').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02415'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '4.45%'
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05216'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '5.41%'
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.006254'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-4.35%'
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007775'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '1.26%'
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1333'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '4.79%'
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007546'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2.31%'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007852'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '11.01%'
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3221'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '10.86%'
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006784'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '5.50%'
$ws.Range('E46').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.08%'
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').Value = 'BOLO'

$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.04498'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '50.60%'
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = 'CoinbaseStockToken'

$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.004202'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.02%'
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.08%'
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.08%'
$ws.Range('E51').Style = 'Normal'
